$wb = $excel.ActiveWorkbook

# Sheet "展览" (overview/exhibition listing) updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 7433
$wsExhibit.Range("F6").Value = 3986
$wsExhibit.Range("F10").Value = 638

# Sheet "全部类型" (all types listing) updates - mirrors the same events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7433
$wsAll.Range("F8").Value = 3986
$wsAll.Range("F12").Value = 638
